$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.441.22'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '3.723.44'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.39'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.93'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '659.51'
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.438'
$ws.Range('E8').Value = '  +3.12%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.07'
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.722.11'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000320'
$ws.Range('E12').Value = '  +18.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.88'
$ws.Range('E13').Value = '  -1.56%  '
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.94'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').Value = '4.417.88'
$ws.Range('E16').Value = '  +1.32%  '
$ws.Range('D17').Value = '97.049.32'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.04'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').Value = '3.718.66'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.06'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.81'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.510'
$ws.Range('E22').Value = '  -4.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '526.63'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000226'
$ws.Range('E25').Value = '  +10.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('E26').Value = '  -3.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '106.70'
$ws.Range('E27').Value = '  +3.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.192'
$ws.Range('E28').Value = '  +14.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.59'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').Value = '3.919.69'
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.03'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.70'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '642.92'
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.594'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.78'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.166'
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.496'
$ws.Range('E43').Value = '  +13.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.76'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.65'
$ws.Range('E45').Value = '  +4.90%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.972'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0456'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.39'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('E51').Value = '  -0.54%  '
